$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Ranking" -> "Rating" (table header-ish single word, isolated run)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Ranking", $true, $false, $false, $false, $false, $true, 1, $false, "Rating", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge "From the start " / "of a new game..." into a single run. This
#    range spans the old _GoBack bookmark, which gets removed as part of
#    the merge -- matching the source edit where the bookmark moved away
#    from here.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("From the start of a new game, a timer counts down how long the game session lasts", $true, $false, $false, $false, $false, $true, 1, $false, "From the start of a new game, a timer counts down how long the game session lasts", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Re-create the _GoBack bookmark right before the "User enters " run
#    (earlier in the document, in the "Dashboard generation" row).
# ---------------------------------------------------------------------------
$bmRange = $d.Content
$bmRange.Find.Execute("User enters", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange.Collapse(1) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 4) Replace remaining "ranking" -> "rating" occurrences (case sensitive,
#    lower-case only -- "Ranking" was already handled above).
# ---------------------------------------------------------------------------
$r = $d.Content
while ($r.Find.Execute("ranking", $true, $false, $false, $false, $false, $true, 1, $false, "rating", 2)) {
    $r.Collapse(0) | Out-Null
}
